$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Plain text / link / name updates (safe as text without coercion)
$ws.Range('D2').Value = '25.757.61'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.626.49'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '1.650.18'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('D14').Value = '1.849.39'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').Value = '0.0₅7543'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '25.785.22'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +5.48%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('D39').Value = '1.109.64'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('E44').Value = '  -2.21%  '
$ws.Range('D45').Value = '1.773.78'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('E46').Value = '  -13.49%  '
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('E51').Value = '  +2.46%  '

# Numeric-looking price updates: force text storage, then restore default style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.18'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5103'
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2560'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06315'
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.39'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07776'
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.222'
$ws.Range('D13').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5500'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.396'
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '193.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.828'
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.981'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '142.03'
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1256'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.729'
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.48'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.237'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04870'
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.221'
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.162'
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.535'
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.379'
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.8910'
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.534'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5492'
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01545'
$ws.Range('D40').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.547'
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7961'
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '97.14'
$ws.Range('D44').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4430'
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.49'
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.518'
$ws.Range('D51').Style = "Normal"
